# Update "想去人数" (interested-people count) values in column F
# for the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1076
    5  = 4629
    7  = 390
    8  = 1370
    9  = 909
    11 = 1064
    13 = 580
    15 = 5
    16 = 263
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
